$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new assignment row's text first (label, then link text) so the
# shared-string table gets the same ordering as the authored workbook.
$ws.Range("A12").Value = "Assignment_11"
$ws.Range("B12").Value = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_11"

# Turn B12 into a real hyperlink pointing at the Assignment_11 folder.
$ws.Hyperlinks.Add($ws.Range("B12"), "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_11") | Out-Null

# Copy formatting (styles) from row 11 (A11:C11) down to row 12, so the new
# row picks up the same cell styles used by the previous "assignment" rows
# (this also overrides the default blue/underlined "Hyperlink" style that
# Hyperlinks.Add just applied to B12).
$ws.Range("A11:C11").Copy() | Out-Null
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Match row 12's height to the other data rows (32.25pt, custom height).
$ws.Rows.Item(12).RowHeight = 32.25

# Set the assignment date for the new row.
$d = Get-Date -Year 2023 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Range("C12").Value = $d

# Update the sheet's current selection to B12.
$ws.Range("B12").Select() | Out-Null
